# Update "想去人数" (number of people interested) counts by +1
# for a set of events that appear on both the "展览" sheet and the
# aggregated "全部类型" sheet.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F5").Value = 80
$wsExhibit.Range("F7").Value = 601
$wsExhibit.Range("F15").Value = 49
$wsExhibit.Range("F18").Value = 296
$wsExhibit.Range("F21").Value = 1104

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F6").Value = 80
$wsAll.Range("F9").Value = 601
$wsAll.Range("F17").Value = 49
$wsAll.Range("F20").Value = 296
$wsAll.Range("F23").Value = 1104
